$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Training Dashboard sheet: update "PERIOD TO EXPIRE" (col H) and
#    "LAST UPDATE" (col I) for rows 3..19 -> progress as of 04-Nov-2025
# ---------------------------------------------------------------------------
$wsTraining = $wb.Worksheets.Item("Training Dashboard")

for ($r = 3; $r -le 19; $r++) {
    $hCell = $wsTraining.Range("H$r")
    $hCell.Value = $hCell.Value2 - 1

    # Column I holds a text-formatted date string, not a real date; force text
    # format so Excel doesn't auto-convert it into a date serial number.
    $iCell = $wsTraining.Range("I$r")
    $iCell.NumberFormat = "@"
    $iCell.Value = "04-Nov-2025"
}

# ---------------------------------------------------------------------------
# 2) Exam Dashboard sheet: insert a new exam record row (Cs Hoist) above the
#    TOTAL AVERAGE row, and refresh the total average value.
# ---------------------------------------------------------------------------
$wsExam = $wb.Worksheets.Item("Exam Dashboard")

# Insert a new blank row at row 6 (pushes the old row 6 "TOTAL AVERAGE" to row 7)
$wsExam.Rows.Item(6).Insert()

# Copy formatting from the row above so the new row matches the rest of the table
$wsExam.Range("A5:G5").Copy()
$wsExam.Range("A6:G6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Columns C (date) and D (percentage) hold text-formatted values, not real
# dates/numbers, so force text format before assigning to avoid Excel
# auto-converting the strings into a date serial / numeric percentage.
$wsExam.Range("C6").NumberFormat = "@"
$wsExam.Range("D6").NumberFormat = "@"
$wsExam.Range("D7").NumberFormat = "@"

# Fill in the new exam record
$wsExam.Range("A6").Value = 4
$wsExam.Range("B6").Value = "Cs Hoist"
$wsExam.Range("C6").Value = "30-Oct-2025"
$wsExam.Range("D6").Value = "90.12%"
$wsExam.Range("E6").Value = "VALID"
$wsExam.Range("F6").Value = "Approved Score. date is valid"

# Update the total average value on the (now shifted) TOTAL AVERAGE row
$wsExam.Range("D7").Value = "90.28%"
